$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1488.4286
$ws.Range("J17").Value = 1133.8148
$ws.Range("L17").Value = 3401.4444
$ws.Range("N17").Value = -3737.4444
$ws.Range("H33").Value = 206.76471
$ws.Range("I33").Value = 264
$ws.Range("K33").Value = 264
$ws.Range("M33").Value = -35
$ws.Range("H43").Value = 1605.5714
$ws.Range("J43").Value = 1673.3334
$ws.Range("L43").Value = 1673.3334
$ws.Range("N43").Value = -1811.3334
$ws.Range("H62").Value = 1865.625
$ws.Range("I62").Value = 1637.2
$ws.Range("J62").Value = 2246.3333
$ws.Range("K62").Value = 1637.2
$ws.Range("L62").Value = 2246.3333
$ws.Range("M62").Value = -1013.2
$ws.Range("N62").Value = -3494.3333
$ws.Range("H65").Value = 1865.625
$ws.Range("I65").Value = 1637.2
$ws.Range("J65").Value = 2246.3333
$ws.Range("K65").Value = 8186
$ws.Range("L65").Value = 11231.6665
$ws.Range("M65").Value = -5066
$ws.Range("N65").Value = -17471.6665
$ws.Range("H96").Value = 1842.8572
$ws.Range("I96").Value = 500
$ws.Range("J96").Value = 3633.3333
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 10899.9999
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -13645.9999
$ws.Range("H106").Value = 4627.25
$ws.Range("I106").Value = 3001
$ws.Range("J106").Value = 9506
$ws.Range("K106").Value = 3001
$ws.Range("L106").Value = 9506
$ws.Range("M106").Value = -2370
$ws.Range("N106").Value = -10768
$ws.Range("H116").Value = 18023.375
$ws.Range("I116").Value = 100000
$ws.Range("K116").Value = 100000
$ws.Range("M116").Value = -96558
$ws.Range("H137").Value = 1714.3846
$ws.Range("I137").Value = 1497.4
$ws.Range("K137").Value = 4492.200000000001
$ws.Range("M137").Value = -1942.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3593.0205
$ws.Range("I32").Value = 2182.4595
$ws.Range("J32").Value = 7942.25
$ws.Range("K32").Value = 2182.4595
$ws.Range("L32").Value = 7942.25
$ws.Range("M32").Value = -1895.4595
$ws.Range("N32").Value = -8516.25
$ws.Range("H45").Value = 2557.7144
$ws.Range("I45").Value = 751
$ws.Range("K45").Value = 751
$ws.Range("M45").Value = -374
$ws.Range("H74").Value = 3867
$ws.Range("I74").Value = 3782.7058
$ws.Range("J74").Value = 4225.25
$ws.Range("K74").Value = 3782.7058
$ws.Range("L74").Value = 4225.25
$ws.Range("M74").Value = -2908.7058
$ws.Range("N74").Value = -5973.25
$ws.Range("H77").Value = 3867
$ws.Range("I77").Value = 3782.7058
$ws.Range("J77").Value = 4225.25
$ws.Range("K77").Value = 18913.529
$ws.Range("L77").Value = 21126.25
$ws.Range("M77").Value = -14545.529
$ws.Range("N77").Value = -29862.25
$ws.Range("H122").Value = 1727.2667
$ws.Range("I122").Value = 1716.3572
$ws.Range("J122").Value = 1880
$ws.Range("K122").Value = 5149.071599999999
$ws.Range("L122").Value = 5640
$ws.Range("M122").Value = -2699.071599999999
$ws.Range("N122").Value = -10540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1409.9048
$ws.Range("J20").Value = 1353.1428
$ws.Range("L20").Value = 1353.1428
$ws.Range("N20").Value = -1847.1428
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 45000
$ws.Range("J55").Value = 45000
$ws.Range("L55").Value = 45000
$ws.Range("N55").Value = -45546
$ws.Range("H105").Value = 2448.5
$ws.Range("I105").Value = 2448.5
$ws.Range("K105").Value = 2448.5
$ws.Range("M105").Value = -701.5
$ws.Range("H134").Value = 6484.4146
$ws.Range("I134").Value = 6843.4062
$ws.Range("K134").Value = 20530.2186
$ws.Range("M134").Value = -17995.2186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 799.75
$ws.Range("I22").Value = 399.66666
$ws.Range("K22").Value = 399.66666
$ws.Range("M22").Value = -49.66665999999998
$ws.Range("H31").Value = 2043.7059
$ws.Range("I31").Value = 981.1539
$ws.Range("J31").Value = 5497
$ws.Range("K31").Value = 981.1539
$ws.Range("L31").Value = 5497
$ws.Range("M31").Value = -686.1539
$ws.Range("N31").Value = -6087
$ws.Range("H34").Value = 2043.7059
$ws.Range("I34").Value = 981.1539
$ws.Range("J34").Value = 5497
$ws.Range("K34").Value = 981.1539
$ws.Range("L34").Value = 5497
$ws.Range("M34").Value = -779.1539
$ws.Range("N34").Value = -5901
$ws.Range("H58").Value = 1307.25
$ws.Range("J58").Value = 1700
$ws.Range("L58").Value = 1700
$ws.Range("N58").Value = -2106
$ws.Range("H132").Value = 2197.125
$ws.Range("I132").Value = 1164.1
$ws.Range("J132").Value = 3918.8333
$ws.Range("K132").Value = 3492.3
$ws.Range("L132").Value = 11756.4999
$ws.Range("M132").Value = -962.2999999999997
$ws.Range("N132").Value = -16816.4999
$ws.Range("H136").Value = 1307.25
$ws.Range("J136").Value = 1700
$ws.Range("L136").Value = 5100
$ws.Range("N136").Value = -10200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H98").Value = 599.25
$ws.Range("I98").Value = 550
$ws.Range("K98").Value = 1650
$ws.Range("M98").Value = -152
$ws.Range("H139").Value = 7364.5293
$ws.Range("I139").Value = 8013.2666
$ws.Range("K139").Value = 24039.7998
$ws.Range("M139").Value = -18899.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4792.8096
$ws.Range("I7").Value = 2751.4285
$ws.Range("J7").Value = 5813.5
$ws.Range("K7").Value = 2751.4285
$ws.Range("L7").Value = 5813.5
$ws.Range("M7").Value = -2639.4285
$ws.Range("N7").Value = -6037.5
$ws.Range("H22").Value = 1068
$ws.Range("I22").Value = 777.2
$ws.Range("K22").Value = 777.2
$ws.Range("M22").Value = -482.2
$ws.Range("H27").Value = 1068
$ws.Range("I27").Value = 777.2
$ws.Range("K27").Value = 777.2
$ws.Range("M27").Value = -670.2
$ws.Range("H126").Value = 4792.8096
$ws.Range("I126").Value = 2751.4285
$ws.Range("J126").Value = 5813.5
$ws.Range("K126").Value = 8254.2855
$ws.Range("L126").Value = 17440.5
$ws.Range("M126").Value = -5784.2855
$ws.Range("N126").Value = -22380.5
$ws.Range("H132").Value = 2014.409
$ws.Range("I132").Value = 1542.7142
$ws.Range("K132").Value = 4628.142599999999
$ws.Range("M132").Value = -2098.142599999999
$ws.Range("H136").Value = 3497.2964
$ws.Range("I136").Value = 2782.238
$ws.Range("K136").Value = 8346.714
$ws.Range("M136").Value = -5796.714
